$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.7172284644194756
$ws1.Range("C2").Value = 0.8431952662721893
$ws1.Range("D2").Value = 0.5337078651685393
$ws1.Range("E2").Value = 0.6536697247706422
$ws1.Range("F2").Value = 0.5759902991107518
$ws1.Range("G2").Value = 0.5413500876680304
$ws1.Range("H2").Value = 0.7172284644194755
$ws1.Range("I2").Value = 285
$ws1.Range("J2").Value = 53
$ws1.Range("K2").Value = 481
$ws1.Range("L2").Value = 249

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 ("0")
$ws2.Range("B2").Value = 0.6589041095890411
$ws2.Range("C2").Value = 0.900749063670412
$ws2.Range("D2").Value = 0.7610759493670886

# row 3 ("1")
$ws2.Range("B3").Value = 0.8431952662721893
$ws2.Range("C3").Value = 0.5337078651685393
$ws2.Range("D3").Value = 0.6536697247706422

# row 4 ("accuracy")
$ws2.Range("B4").Value = 0.7172284644194756
$ws2.Range("C4").Value = 0.7172284644194756
$ws2.Range("D4").Value = 0.7172284644194756
$ws2.Range("E4").Value = 0.7172284644194756

# row 5 ("macro avg")
$ws2.Range("B5").Value = 0.7510496879306152
$ws2.Range("C5").Value = 0.7172284644194756
$ws2.Range("D5").Value = 0.7073728370688654

# row 6 ("weighted avg")
$ws2.Range("B6").Value = 0.7510496879306152
$ws2.Range("C6").Value = 0.7172284644194756
$ws2.Range("D6").Value = 0.7073728370688653

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# row 2 ("Predicted 0")
$ws3.Range("B2").Value = 481
$ws3.Range("C2").Value = 53

# row 3 ("Predicted 1")
$ws3.Range("B3").Value = 249
$ws3.Range("C3").Value = 285
